$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 71: homework entry for 2020-02-27 (ENEST / 03025)
$ws.Range("A71").Value = 1582761600
$ws.Range("B71").Value = "'2020-02-27"
$ws.Range("C71").Value = "'03025"
$ws.Range("D71").Value = "ENEST"
$ws.Range("E71").Value = 0.13
$ws.Range("F71").Value = 0.13
$ws.Range("G71").Value = 0.13
$ws.Range("H71").Value = 0.13
$ws.Range("I71").Value = "-"
